# "Uren Registratie en prezi"
# Fill in the hours that were logged on the Friday of Week 16 (row 94) and
# the Monday of Week 17 (row 98), plus the week-17 manual weekly total
# (B103). The week-16 manual weekly total (B95) is bumped too, since Carlo
# apparently logged extra hours that week. All the SUM()-formula cells
# (C95:I95, C103:I103) and the summary block in rows 2-10 (K2, L2:L8,
# M2:M10) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 16 ("Vrijdag", row 94): Zinedine, Robin, Carlo worked ---
$ws.Range("D94").Value = 3
$ws.Range("E94").Value = 4
$ws.Range("F94").Value = 4

# --- Week 16 weekly manual total (column B is entered by hand, not summed) ---
$ws.Range("B95").Value = 16

# --- Week 17 ("Maandag", row 98): Rief, Zinedine, Robin, Carlo, Sam worked ---
$ws.Range("C98").Value = 4
$ws.Range("D98").Value = 4
$ws.Range("E98").Value = 4
$ws.Range("F98").Value = 4
$ws.Range("H98").Value = 3

# --- Week 17 weekly manual total (column B is entered by hand, not summed) ---
$ws.Range("B103").Value = 4

# --- Scroll / selection bookkeeping (matches the saved view state) ---
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K92").Select()
